$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows starting at row 3. This pushes the existing row 3
# (www.als.com_6.json) down to row 6, while row 2 (soxy.com.json) stays at row 2
# for now (it will be moved to row 5 below).
$ws.Rows("3:5").Insert()

# Row 2: replace soxy.com.json data with the new first entry
$ws.Range("A2").Value = "1999beauty.com.json"
$ws.Range("B2").Value = 26
$ws.Range("C2").Value = 2

# Row 3: new entry
$ws.Range("A3").Value = "786cosmetics.com_2.json"
$ws.Range("B3").Value = 94
$ws.Range("C3").Value = 88

# Row 4: new entry
$ws.Range("A4").Value = "beautysociety.com_9.json"
$ws.Range("B4").Value = 17
$ws.Range("C4").Value = 11

# Row 5: the original soxy.com.json data, now moved down
$ws.Range("A5").Value = "soxy.com.json"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0

# Row 6 already holds www.als.com_6.json / 877 / 656 thanks to the insert shifting it down.

Write-Output "done"
